$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"32.50235"
$ws.Range("H2").Value = [double]"97.50704999999999"
$ws.Range("I2").Value = [double]"0.004318312013857221"
$ws.Range("J2").Value = [double]"0.004318312013857221"
$ws.Range("M2").Value = [double]"101.1420973333333"
$ws.Range("N2").Value = [double]"303.426292"
$ws.Range("O2").Value = [double]"0.6720924517779291"
$ws.Range("P2").Value = [double]"0.6720924517779291"
$ws.Range("Q2").Value = [double]"3287.355847262066"
$ws.Range("R2").Value = [double]"29586.2026253586"
$ws.Range("S2").Value = [double]"0.002902304908935386"
$ws.Range("T2").Value = [double]"0.002902304908935386"
$ws.Range("G3").Value = [double]"32.50235"
$ws.Range("H3").Value = [double]"97.50704999999999"
$ws.Range("I3").Value = [double]"0.004318312013857221"
$ws.Range("J3").Value = [double]"0.004318312013857221"
$ws.Range("O3").Value = [double]"0.007736938980150513"
$ws.Range("P3").Value = [double]"0.007736938980150512"
$ws.Range("Q3").Value = [double]"37.84311448376666"
$ws.Range("R3").Value = [double]"340.5880303539"
$ws.Range("S3").Value = [double]"3.34105165484642E-05"
$ws.Range("T3").Value = [double]"3.341051654846419E-05"
$ws.Range("G4").Value = [double]"32.50235"
$ws.Range("H4").Value = [double]"97.50704999999999"
$ws.Range("I4").Value = [double]"0.004318312013857221"
$ws.Range("J4").Value = [double]"0.004318312013857221"
$ws.Range("M4").Value = [double]"26.91044433333333"
$ws.Range("N4").Value = [double]"80.73133300000001"
$ws.Range("O4").Value = [double]"0.1788207579957193"
$ws.Range("P4").Value = [double]"0.1788207579957193"
$ws.Range("Q4").Value = [double]"874.6526803775167"
$ws.Range("R4").Value = [double]"7871.87412339765"
$ws.Range("S4").Value = [double]"0.0007722038275799692"
$ws.Range("T4").Value = [double]"0.0007722038275799692"
$ws.Range("G5").Value = [double]"32.50235"
$ws.Range("H5").Value = [double]"97.50704999999999"
$ws.Range("I5").Value = [double]"0.004318312013857221"
$ws.Range("J5").Value = [double]"0.004318312013857221"
$ws.Range("M5").Value = [double]"0.7853516666666667"
$ws.Range("N5").Value = [double]"2.356055"
$ws.Range("O5").Value = [double]"0.00521868678892747"
$ws.Range("P5").Value = [double]"0.00521868678892747"
$ws.Range("Q5").Value = [double]"25.52577474308333"
$ws.Range("R5").Value = [double]"229.73197268775"
$ws.Range("S5").Value = [double]"2.253591785718346E-05"
$ws.Range("T5").Value = [double]"2.253591785718346E-05"
$ws.Range("G6").Value = [double]"32.50235"
$ws.Range("H6").Value = [double]"97.50704999999999"
$ws.Range("I6").Value = [double]"0.004318312013857221"
$ws.Range("J6").Value = [double]"0.004318312013857221"
$ws.Range("M6").Value = [double]"20.486157"
$ws.Range("N6").Value = [double]"61.458471"
$ws.Range("O6").Value = [double]"0.1361311644572737"
$ws.Range("P6").Value = [double]"0.1361311644572737"
$ws.Range("Q6").Value = [double]"665.8482449689501"
$ws.Range("R6").Value = [double]"5992.63420472055"
$ws.Range("S6").Value = [double]"0.0005878568429362183"
$ws.Range("T6").Value = [double]"0.0005878568429362182"
$ws.Range("I7").Value = [double]"0.006762540683959845"
$ws.Range("J7").Value = [double]"0.006762540683959845"
$ws.Range("M7").Value = [double]"101.1420973333333"
$ws.Range("N7").Value = [double]"303.426292"
$ws.Range("O7").Value = [double]"0.6720924517779291"
$ws.Range("P7").Value = [double]"0.6720924517779291"
$ws.Range("Q7").Value = [double]"5148.048030903133"
$ws.Range("R7").Value = [double]"46332.4322781282"
$ws.Range("S7").Value = [double]"0.004545052548530566"
$ws.Range("T7").Value = [double]"0.004545052548530566"
$ws.Range("I8").Value = [double]"0.006762540683959845"
$ws.Range("J8").Value = [double]"0.006762540683959845"
$ws.Range("O8").Value = [double]"0.007736938980150513"
$ws.Range("P8").Value = [double]"0.007736938980150512"
$ws.Range("S8").Value = [double]"5.232136462258264E-05"
$ws.Range("T8").Value = [double]"5.232136462258264E-05"
$ws.Range("I9").Value = [double]"0.006762540683959845"
$ws.Range("J9").Value = [double]"0.006762540683959845"
$ws.Range("M9").Value = [double]"26.91044433333333"
$ws.Range("N9").Value = [double]"80.73133300000001"
$ws.Range("O9").Value = [double]"0.1788207579957193"
$ws.Range("P9").Value = [double]"0.1788207579957193"
$ws.Range("Q9").Value = [double]"1369.719074584463"
$ws.Range("R9").Value = [double]"12327.47167126017"
$ws.Range("S9").Value = [double]"0.001209282651082589"
$ws.Range("T9").Value = [double]"0.001209282651082589"
$ws.Range("I10").Value = [double]"0.006762540683959845"
$ws.Range("J10").Value = [double]"0.006762540683959845"
$ws.Range("M10").Value = [double]"0.7853516666666667"
$ws.Range("N10").Value = [double]"2.356055"
$ws.Range("O10").Value = [double]"0.00521868678892747"
$ws.Range("P10").Value = [double]"0.00521868678892747"
$ws.Range("Q10").Value = [double]"39.97374197042055"
$ws.Range("R10").Value = [double]"359.763677733785"
$ws.Range("S10").Value = [double]"3.529158172696578E-05"
$ws.Range("T10").Value = [double]"3.529158172696578E-05"
$ws.Range("I11").Value = [double]"0.006762540683959845"
$ws.Range("J11").Value = [double]"0.006762540683959845"
$ws.Range("M11").Value = [double]"20.486157"
$ws.Range("N11").Value = [double]"61.458471"
$ws.Range("O11").Value = [double]"0.1361311644572737"
$ws.Range("P11").Value = [double]"0.1361311644572737"
$ws.Range("Q11").Value = [double]"1042.728230729153"
$ws.Range("R11").Value = [double]"9384.554076562377"
$ws.Range("S11").Value = [double]"0.0009205925379971422"
$ws.Range("T11").Value = [double]"0.000920592537997142"
$ws.Range("G12").Value = [double]"3274.382486666667"
$ws.Range("H12").Value = [double]"9823.14746"
$ws.Range("I12").Value = [double]"0.4350394734576531"
$ws.Range("J12").Value = [double]"0.435039473457653"
$ws.Range("M12").Value = [double]"101.1420973333333"
$ws.Range("N12").Value = [double]"303.426292"
$ws.Range("O12").Value = [double]"0.6720924517779291"
$ws.Range("P12").Value = [double]"0.6720924517779291"
$ws.Range("Q12").Value = [double]"331177.912173002"
$ws.Range("R12").Value = [double]"2980601.209557018"
$ws.Range("S12").Value = [double]"0.2923867463363334"
$ws.Range("T12").Value = [double]"0.2923867463363333"
$ws.Range("G13").Value = [double]"3274.382486666667"
$ws.Range("H13").Value = [double]"9823.14746"
$ws.Range("I13").Value = [double]"0.4350394734576531"
$ws.Range("J13").Value = [double]"0.435039473457653"
$ws.Range("O13").Value = [double]"0.007736938980150513"
$ws.Range("P13").Value = [double]"0.007736938980150512"
$ws.Range("Q13").Value = [double]"3812.426833954075"
$ws.Range("R13").Value = [double]"34311.84150558668"
$ws.Range("S13").Value = [double]"0.003365873860098671"
$ws.Range("T13").Value = [double]"0.00336587386009867"
$ws.Range("G14").Value = [double]"3274.382486666667"
$ws.Range("H14").Value = [double]"9823.14746"
$ws.Range("I14").Value = [double]"0.4350394734576531"
$ws.Range("J14").Value = [double]"0.435039473457653"
$ws.Range("M14").Value = [double]"26.91044433333333"
$ws.Range("N14").Value = [double]"80.73133300000001"
$ws.Range("O14").Value = [double]"0.1788207579957193"
$ws.Range("P14").Value = [double]"0.1788207579957193"
$ws.Range("Q14").Value = [double]"88115.08763348492"
$ws.Range("R14").Value = [double]"793035.7887013643"
$ws.Range("S14").Value = [double]"0.07779408840175611"
$ws.Range("T14").Value = [double]"0.0777940884017561"
$ws.Range("G15").Value = [double]"3274.382486666667"
$ws.Range("H15").Value = [double]"9823.14746"
$ws.Range("I15").Value = [double]"0.4350394734576531"
$ws.Range("J15").Value = [double]"0.435039473457653"
$ws.Range("M15").Value = [double]"0.7853516666666667"
$ws.Range("N15").Value = [double]"2.356055"
$ws.Range("O15").Value = [double]"0.00521868678892747"
$ws.Range("P15").Value = [double]"0.00521868678892747"
$ws.Range("Q15").Value = [double]"2571.541743207811"
$ws.Range("R15").Value = [double]"23143.8756888703"
$ws.Range("S15").Value = [double]"0.002270334752795417"
$ws.Range("T15").Value = [double]"0.002270334752795416"
$ws.Range("G16").Value = [double]"3274.382486666667"
$ws.Range("H16").Value = [double]"9823.14746"
$ws.Range("I16").Value = [double]"0.4350394734576531"
$ws.Range("J16").Value = [double]"0.435039473457653"
$ws.Range("M16").Value = [double]"20.486157"
$ws.Range("N16").Value = [double]"61.458471"
$ws.Range("O16").Value = [double]"0.1361311644572737"
$ws.Range("P16").Value = [double]"0.1361311644572737"
$ws.Range("Q16").Value = [double]"67079.51369990375"
$ws.Range("R16").Value = [double]"603715.6232991337"
$ws.Range("S16").Value = [double]"0.05922243010666955"
$ws.Range("T16").Value = [double]"0.05922243010666953"
$ws.Range("G17").Value = [double]"7.278837333333333"
$ws.Range("H17").Value = [double]"21.836512"
$ws.Range("I17").Value = [double]"0.0009670774791190726"
$ws.Range("J17").Value = [double]"0.0009670774791190726"
$ws.Range("M17").Value = [double]"101.1420973333333"
$ws.Range("N17").Value = [double]"303.426292"
$ws.Range("O17").Value = [double]"0.6720924517779291"
$ws.Range("P17").Value = [double]"0.6720924517779291"
$ws.Range("Q17").Value = [double]"736.1968740415003"
$ws.Range("R17").Value = [double]"6625.771866373503"
$ws.Range("S17").Value = [double]"0.0006499654740003566"
$ws.Range("T17").Value = [double]"0.0006499654740003566"
$ws.Range("G18").Value = [double]"7.278837333333333"
$ws.Range("H18").Value = [double]"21.836512"
$ws.Range("I18").Value = [double]"0.0009670774791190726"
$ws.Range("J18").Value = [double]"0.0009670774791190726"
$ws.Range("O18").Value = [double]"0.007736938980150513"
$ws.Range("P18").Value = [double]"0.007736938980150512"
$ws.Range("Q18").Value = [double]"8.474891031388443"
$ws.Range("R18").Value = [double]"76.27401928249598"
$ws.Range("S18").Value = [double]"7.482219445022048E-06"
$ws.Range("T18").Value = [double]"7.482219445022047E-06"
$ws.Range("G19").Value = [double]"7.278837333333333"
$ws.Range("H19").Value = [double]"21.836512"
$ws.Range("I19").Value = [double]"0.0009670774791190726"
$ws.Range("J19").Value = [double]"0.0009670774791190726"
$ws.Range("M19").Value = [double]"26.91044433333333"
$ws.Range("N19").Value = [double]"80.73133300000001"
$ws.Range("O19").Value = [double]"0.1788207579957193"
$ws.Range("P19").Value = [double]"0.1788207579957193"
$ws.Range("Q19").Value = [double]"195.8767468700551"
$ws.Range("R19").Value = [double]"1762.890721830496"
$ws.Range("S19").Value = [double]"0.0001729335278566619"
$ws.Range("T19").Value = [double]"0.0001729335278566619"
$ws.Range("G20").Value = [double]"7.278837333333333"
$ws.Range("H20").Value = [double]"21.836512"
$ws.Range("I20").Value = [double]"0.0009670774791190726"
$ws.Range("J20").Value = [double]"0.0009670774791190726"
$ws.Range("M20").Value = [double]"0.7853516666666667"
$ws.Range("N20").Value = [double]"2.356055"
$ws.Range("O20").Value = [double]"0.00521868678892747"
$ws.Range("P20").Value = [double]"0.00521868678892747"
$ws.Range("Q20").Value = [double]"5.716447031128888"
$ws.Range("R20").Value = [double]"51.44802328016"
$ws.Range("S20").Value = [double]"5.046874464147986E-06"
$ws.Range("T20").Value = [double]"5.046874464147986E-06"
$ws.Range("G21").Value = [double]"7.278837333333333"
$ws.Range("H21").Value = [double]"21.836512"
$ws.Range("I21").Value = [double]"0.0009670774791190726"
$ws.Range("J21").Value = [double]"0.0009670774791190726"
$ws.Range("M21").Value = [double]"20.486157"
$ws.Range("N21").Value = [double]"61.458471"
$ws.Range("O21").Value = [double]"0.1361311644572737"
$ws.Range("P21").Value = [double]"0.1361311644572737"
$ws.Range("Q21").Value = [double]"149.115404388128"
$ws.Range("R21").Value = [double]"1342.038639493152"
$ws.Range("S21").Value = [double]"0.0001316493833528842"
$ws.Range("T21").Value = [double]"0.0001316493833528842"
$ws.Range("G22").Value = [double]"4161.570231333333"
$ws.Range("H22").Value = [double]"12484.710694"
$ws.Range("I22").Value = [double]"0.5529125963654108"
$ws.Range("J22").Value = [double]"0.5529125963654108"
$ws.Range("M22").Value = [double]"101.1420973333333"
$ws.Range("N22").Value = [double]"303.426292"
$ws.Range("O22").Value = [double]"0.6720924517779291"
$ws.Range("P22").Value = [double]"0.6720924517779291"
$ws.Range("Q22").Value = [double]"420909.9413970184"
$ws.Range("R22").Value = [double]"3788189.472573166"
$ws.Range("S22").Value = [double]"0.3716083825101294"
$ws.Range("T22").Value = [double]"0.3716083825101294"
$ws.Range("G23").Value = [double]"4161.570231333333"
$ws.Range("H23").Value = [double]"12484.710694"
$ws.Range("I23").Value = [double]"0.5529125963654108"
$ws.Range("J23").Value = [double]"0.5529125963654108"
$ws.Range("O23").Value = [double]"0.007736938980150513"
$ws.Range("P23").Value = [double]"0.007736938980150512"
$ws.Range("Q23").Value = [double]"4845.396677365872"
$ws.Range("R23").Value = [double]"43608.57009629285"
$ws.Range("S23").Value = [double]"0.004277851019435774"
$ws.Range("T23").Value = [double]"0.004277851019435773"
$ws.Range("G24").Value = [double]"4161.570231333333"
$ws.Range("H24").Value = [double]"12484.710694"
$ws.Range("I24").Value = [double]"0.5529125963654108"
$ws.Range("J24").Value = [double]"0.5529125963654108"
$ws.Range("M24").Value = [double]"26.91044433333333"
$ws.Range("N24").Value = [double]"80.73133300000001"
$ws.Range("O24").Value = [double]"0.1788207579957193"
$ws.Range("P24").Value = [double]"0.1788207579957193"
$ws.Range("Q24").Value = [double]"111989.7040495528"
$ws.Range("R24").Value = [double]"1007907.336445975"
$ws.Range("S24").Value = [double]"0.09887224958744392"
$ws.Range("T24").Value = [double]"0.09887224958744392"
$ws.Range("G25").Value = [double]"4161.570231333333"
$ws.Range("H25").Value = [double]"12484.710694"
$ws.Range("I25").Value = [double]"0.5529125963654108"
$ws.Range("J25").Value = [double]"0.5529125963654108"
$ws.Range("M25").Value = [double]"0.7853516666666667"
$ws.Range("N25").Value = [double]"2.356055"
$ws.Range("O25").Value = [double]"0.00521868678892747"
$ws.Range("P25").Value = [double]"0.00521868678892747"
$ws.Range("Q25").Value = [double]"3268.296117128019"
$ws.Range("R25").Value = [double]"29414.66505415217"
$ws.Range("S25").Value = [double]"0.002885477662083756"
$ws.Range("T25").Value = [double]"0.002885477662083756"
$ws.Range("G26").Value = [double]"4161.570231333333"
$ws.Range("H26").Value = [double]"12484.710694"
$ws.Range("I26").Value = [double]"0.5529125963654108"
$ws.Range("J26").Value = [double]"0.5529125963654108"
$ws.Range("M26").Value = [double]"20.486157"
$ws.Range("N26").Value = [double]"61.458471"
$ws.Range("O26").Value = [double]"0.1361311644572737"
$ws.Range("P26").Value = [double]"0.1361311644572737"
$ws.Range("Q26").Value = [double]"85254.58112562099"
$ws.Range("R26").Value = [double]"767291.2301305889"
$ws.Range("S26").Value = [double]"0.07526863558631795"
$ws.Range("T26").Value = [double]"0.07526863558631794"
